# Apply the refreshed crypto Price (column D) and Volume(1h) change
# percentage (column E) figures to Sheet1, as produced by the scheduled
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe forces Excel to store the Price column values as
# literal text (matching the original plain-text price cells in the sheet)
# instead of re-parsing number-look-alikes (e.g. "3.085.49" or "1.00")
# into floating point values and losing formatting/precision.

$ws.Range('D2').Value = "'" + "63.487.77"
$ws.Range('E2').Value = "  +0.03%  "
$ws.Range('D3').Value = "'" + "3.085.49"
$ws.Range('E3').Value = "  -0.20%  "
$ws.Range('E4').Value = "  -0.06%  "
$ws.Range('D5').Value = "'" + "547.57"
$ws.Range('E5').Value = "  -0.26%  "
$ws.Range('D6').Value = "'" + "139.64"
$ws.Range('E6').Value = "  +1.73%  "
$ws.Range('E7').Value = "  -0.08%  "
$ws.Range('D8').Value = "'" + "3.079.17"
$ws.Range('E8').Value = "  -0.36%  "
$ws.Range('D9').Value = "'" + "0.498"
$ws.Range('E9').Value = "  +0.28%  "
$ws.Range('E10').Value = "  +1.20%  "
$ws.Range('D11').Value = "'" + "6.42"
$ws.Range('E11').Value = "  +1.23%  "
$ws.Range('E12').Value = "  -2.73%  "
$ws.Range('E13').Value = "  +3.76%  "
$ws.Range('D14').Value = "'" + "35.06"
$ws.Range('E14').Value = "  -1.19%  "
$ws.Range('D15').Value = "'" + "3.586.08"
$ws.Range('E15').Value = "  -0.35%  "
$ws.Range('D16').Value = "'" + "63.517.41"
$ws.Range('E16').Value = "  +0.09%  "
$ws.Range('E17').Value = "  +1.08%  "
$ws.Range('D18').Value = "'" + "3.082.78"
$ws.Range('E18').Value = "  -0.26%  "
$ws.Range('E19').Value = "  -1.35%  "
$ws.Range('D20').Value = "'" + "475.86"
$ws.Range('E20').Value = "  -2.70%  "
$ws.Range('D21').Value = "'" + "13.53"
$ws.Range('E21').Value = "  -0.64%  "
$ws.Range('D22').Value = "'" + "0.704"
$ws.Range('E22').Value = "  -2.04%  "
$ws.Range('E23').Value = "  -2.63%  "
$ws.Range('E24').Value = "  -0.43%  "
$ws.Range('D25').Value = "'" + "12.26"
$ws.Range('E25').Value = "  -0.93%  "
$ws.Range('E26').Value = "  +0.05%  "
$ws.Range('E27').Value = "  -1.31%  "
$ws.Range('E28').Value = "  -6.44%  "
$ws.Range('D29').Value = "'" + "1.00"
$ws.Range('E29').Value = "  -0.03%  "
$ws.Range('D30').Value = "'" + "26.31"
$ws.Range('E30').Value = "  -1.33%  "
$ws.Range('E31').Value = "  -3.48%  "
$ws.Range('E32').Value = "  +2.95%  "
$ws.Range('D33').Value = "'" + "59.17"
$ws.Range('E33').Value = "  +0.25%  "
$ws.Range('E34').Value = "  -7.34%  "
$ws.Range('D35').Value = "'" + "5.54"
$ws.Range('E35').Value = "  +8.26%  "
$ws.Range('D36').Value = "'" + "6.04"
$ws.Range('E36').Value = "  -0.47%  "
$ws.Range('D37').Value = "'" + "490.24"
$ws.Range('E37').Value = "  -3.33%  "
$ws.Range('D38').Value = "'" + "3.264.83"
$ws.Range('E38').Value = "  +3.71%  "
$ws.Range('E39').Value = "  +0.96%  "
$ws.Range('E40').Value = "  -0.44%  "
$ws.Range('E41').Value = "  -0.83%  "
$ws.Range('D42').Value = "'" + "8.17"
$ws.Range('E42').Value = "  -0.01%  "
$ws.Range('D43').Value = "'" + "2.61"
$ws.Range('E43').Value = "  -0.85%  "
$ws.Range('E44').Value = "  -1.11%  "
$ws.Range('E45').Value = "  +0.01%  "
$ws.Range('D46').Value = "'" + "25.51"
$ws.Range('E46').Value = "  +0.41%  "
$ws.Range('D47').Value = "'" + "124.55"
$ws.Range('E47').Value = "  +3.46%  "
$ws.Range('E48').Value = "  -1.69%  "
$ws.Range('E49').Value = "  +4.89%  "
$ws.Range('D50').Value = "'" + "0.109"
$ws.Range('E50').Value = "  +0.38%  "
$ws.Range('E51').Value = "  -0.17%  "
